$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the team-record columns (AD, AE, AF).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing bold/bordered/centered header style (same as the rest
# of row 1, e.g. AC1) by copy/paste-special-formats rather than touching
# ".Style" directly (which doesn't stick in this runtime).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team's W/L/T record for every player row (2-59).
$lastRow = 59
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 88
    $ws.Cells.Item($r, 31).Value = 74
    $ws.Cells.Item($r, 32).Value = 0
}
